# Apply the "cn181107" console upload edit to elasticDetailTabsafetyDroup.xlsx
#
# Summary of the change (from the OOXML diff):
#  1. The old (grammatically awkward) English translation for the "security
#     group cap reached" string is replaced by a cleaner sentence. Because the
#     shared-string table is de-duplicated/compacted on save, simply writing
#     the new text into C1 removes the stale <si> and appends a fresh one,
#     which automatically renumbers every other C-column shared-string index
#     exactly like the diff shows (12->17, 13->12, 14->13, 15->14, 16->15,
#     17->16).
#  2. Column B is widened slightly (custom width) while column C keeps the
#     sheet's default width.
#  3. The selected cell moves from C5 to C15, and the view is scrolled so
#     column B is the first visible column.
#  4. A page setup (A4, portrait) is defined for the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the English translation in C1 (row for i18nKey_1 / cap-reached message)
$ws.Range("C1").Value = "Unable to add security groups since the upper limit has been reached."

# 2. Widen column B; leave column C at the sheet default width.
$ws.Columns.Item(2).ColumnWidth = 31.5

# 3. Update the view: scroll so column B is the left-most visible column,
#    then select C15.
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C15").Select()

# 4. Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
